$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name to reflect new date
$ws.Name = "Through 2022-06-13"

# Update the "June (through ...)" label in A7
$ws.Range("A7").Value = "June (through 06-13)"

# Update June row (row 7) values
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 17
$ws.Range("D7").Value = 22
$ws.Range("E7").Value = 25
$ws.Range("F7").Value = 17
$ws.Range("G7").Value = 51
$ws.Range("H7").Value = 47
$ws.Range("I7").Value = 59

# Update Total row (row 8) values
$ws.Range("B8").Value = 113
$ws.Range("C8").Value = 226
$ws.Range("D8").Value = 338
$ws.Range("E8").Value = 320
$ws.Range("F8").Value = 221
$ws.Range("G8").Value = 409
$ws.Range("H8").Value = 678
$ws.Range("I8").Value = 722
